$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Authorship value (name change: Villars -> Villars-Amberg)
$ws.Range("H2").Value2 = "Daniela Subotic, Noémi Villars-Amberg"

# Add a new "Authorship Resource" column (I), matching the header style of
# the existing "Authorship" header (H1) which has a bottom border + bold font.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "Authorship Resource"
$ws.Range("I2").Value2 = "Daniela Subotic, Noémi Villars-Amberg"

# Size the new column similarly to the other descriptive columns.
$ws.Columns.Item(9).ColumnWidth = 44.66666666666667

# Clear the stray "general alignment" formatting that used to sit on most of
# row 2 (column D kept its normal bordered style, the rest go back to plain).
$ws.Range("A2:C2").Style = "Normal"
$ws.Range("E2:I2").Style = "Normal"

# Restore the selection that was active when the workbook was last saved.
$ws.Range("C15").Select() | Out-Null
